$p = $ppt.ActivePresentation

# Slide 8 ("Maintenance Report"): replace the single paragraph in the body
# placeholder with a short bulleted breakdown of the changes.
$s = $p.Slides.Item(8)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Text = "Changes`rNaming of the TOE`rAdded a manufacturing site"

# Indent the two detail bullets one level under "Changes".
$tr.Paragraphs(2,1).IndentLevel = 2
$tr.Paragraphs(3,1).IndentLevel = 2
